$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing company entry, add new metric columns ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("F2").Value = 0.301
$ws.Range("G2").Value = 0.07999615056117601
$ws.Range("H2").Value = 0.07999615056117601
$ws.Range("I2").Value = -0.147401428302197
$ws.Range("J2").Value = -0.1259592075101148
$ws.Range("K2").Value = 141.61
$ws.Range("L2").Value = 0.5678323248605581
$ws.Range("M2").Value = 13.6
$ws.Range("N2").Value = 0.0005247350859485204
$ws.Range("O2").Value = 0.09603841536614645
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 13.6
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 390.227
$ws.Range("V2").Value = 0.01505630870473774
$ws.Range("W2").Value = -2.027973258658284
$ws.Range("X2").Value = 0.07145886537644609
$ws.Range("Y2").Value = -2.09943212403473
$ws.Range("Z2").Value = 0.09868747521802204
$ws.Range("AA2").Value = -2.312631366258389
$ws.Range("AB2").Value = 0.0704384594007583
$ws.Range("AC2").Value = -2.383069825659147
$ws.Range("AD2").Value = 1048.423
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1048.423
$ws.Range("AG2").Value = 658.196
$ws.Range("AH2").Value = 0.03887906158891946
$ws.Range("AI2").Value = 0.2875209486118314
$ws.Range("AJ2").Value = 0.02476652274251886
$ws.Range("AK2").Value = 0.2021364807278186
$ws.Range("AL2").Value = 68.10199999999999
$ws.Range("AM2").Value = -207.798
$ws.Range("AN2").Value = -78.65138784696174
$ws.Range("AO2").Value = -0.5397785674429534
$ws.Range("AP2").Value = -49.37704426106527
$ws.Range("AQ2").Value = 0.1769025688408936

# --- Row 3: rename company, update metric columns ---
$ws.Range("B3").Value = "StoneCo Ltd. (NasdaqGS:STNE)"
$ws.Range("F3").Value = 0.301
$ws.Range("G3").Value = 0.08688656476267097
$ws.Range("H3").Value = 0.08688656476267097
$ws.Range("I3").Value = -0.1415929203539823
$ws.Range("J3").Value = -0.1003983851884021
$ws.Range("K3").Value = 143
$ws.Range("L3").Value = 0.5752212389380531
$ws.Range("M3").Value = 13.6
$ws.Range("N3").Value = 0.0005249323570620772
$ws.Range("O3").Value = 0.0951048951048951
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 13.6
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 390.2
$ws.Range("V3").Value = 0.01506092689158989
$ws.Range("W3").Value = 0.1057301293900185
$ws.Range("X3").Value = 0.07230155806225204
$ws.Range("Y3").Value = 0.03342857132776644
$ws.Range("Z3").Value = 0.09838920330866348
$ws.Range("AA3").Value = -0.009878117132163199
$ws.Range("AB3").Value = 0.07070089954247212
$ws.Range("AC3").Value = -0.08057901667463532
$ws.Range("AD3").Value = 1048.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1048.3
$ws.Range("AG3").Value = 658.0999999999999
$ws.Range("AH3").Value = 0.03888872401359232
$ws.Range("AI3").Value = 0.2874102100126117
$ws.Range("AJ3").Value = 0.02477207880690501
$ws.Range("AK3").Value = 0.2020447009701584
$ws.Range("AL3").Value = 68.09999999999999
$ws.Range("AM3").Value = -207.8
$ws.Range("AN3").Value = -88.83898305084745
$ws.Range("AO3").Value = -0.5168869309838474
$ws.Range("AP3").Value = -55.77118644067796
$ws.Range("AQ3").Value = 0.1693936477382098

# --- Row 4: new row for Sealand Capital Galaxy Limited (re-added with updated figures) ---
$ws.Range("A4").Value = "Cayman Islands"
$ws.Range("B4").Value = "Sealand Capital Galaxy Limited (LSE:SCGL)"
$ws.Range("C4").Value = "Information Services"
$ws.Range("G4").Value = -2.096569250317662
$ws.Range("H4").Value = -2.096569250317662
$ws.Range("I4").Value = -1.982210927573062
$ws.Range("J4").Value = -1.982210927573062
$ws.Range("K4").Value = -1.39
$ws.Range("L4").Value = -1.766200762388818
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.027
$ws.Range("V4").Value = 0.002772073921971252
$ws.Range("W4").Value = -4.161676646706586
$ws.Range("X4").Value = 0.07061617269064015
$ws.Range("Y4").Value = -4.232292819397226
$ws.Range("Z4").Value = 2.328402366863905
$ws.Range("AA4").Value = -4.615384615384615
$ws.Range("AB4").Value = 0.07017601925904446
$ws.Range("AC4").Value = -4.685560634643659
$ws.Range("AD4").Value = 0.123
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.123
$ws.Range("AG4").Value = 0.096
$ws.Range("AH4").Value = 0.01247085065395924
$ws.Range("AI4").Value = -0.1258955987717502
$ws.Range("AJ4").Value = 0.009760065067100447
$ws.Range("AK4").Value = -0.09561752988047809
$ws.Range("AL4").Value = 0.002
$ws.Range("AM4").Value = 0.002
$ws.Range("AN4").Value = -0.08039215686274509
$ws.Range("AO4").Value = -780
$ws.Range("AP4").Value = -0.06274509803921569
$ws.Range("AQ4").Value = -780
